$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 17.05724733333333
$ws.Range("H2").Value = 51.17174199999999
$ws.Range("I2").Value = 0.6861755591406274
$ws.Range("J2").Value = 0.6861755591406274
$ws.Range("M2").Value = 8.488196666666667
$ws.Range("N2").Value = 25.46459
$ws.Range("O2").Value = 0.04138402976425696
$ws.Range("P2").Value = 0.04138402976425696
$ws.Range("Q2").Value = 144.7852699573089
$ws.Range("R2").Value = 1303.06742961578
$ws.Range("S2").Value = 0.02839670976298139
$ws.Range("T2").Value = 0.02839670976298139
# Row 3
$ws.Range("G3").Value = 17.05724733333333
$ws.Range("H3").Value = 51.17174199999999
$ws.Range("I3").Value = 0.6861755591406274
$ws.Range("J3").Value = 0.6861755591406274
$ws.Range("O3").Value = 0.3297460182766552
$ws.Range("P3").Value = 0.3297460182766552
$ws.Range("Q3").Value = 1153.642275667603
$ws.Range("R3").Value = 10382.78048100843
$ws.Range("S3").Value = 0.2262636584653794
$ws.Range("T3").Value = 0.2262636584653794
# Row 4
$ws.Range("G4").Value = 17.05724733333333
$ws.Range("H4").Value = 51.17174199999999
$ws.Range("I4").Value = 0.6861755591406274
$ws.Range("J4").Value = 0.6861755591406274
$ws.Range("O4").Value = 0.6288699519590879
$ws.Range("P4").Value = 0.6288699519590879
$ws.Range("Q4").Value = 2200.150789594601
$ws.Range("R4").Value = 19801.35710635141
$ws.Range("S4").Value = 0.4315151909122666
$ws.Range("T4").Value = 0.4315151909122666
# Row 5
$ws.Range("I5").Value = 0.1466018818485066
$ws.Range("J5").Value = 0.1466018818485066
$ws.Range("M5").Value = 8.488196666666667
$ws.Range("N5").Value = 25.46459
$ws.Range("O5").Value = 0.04138402976425696
$ws.Range("P5").Value = 0.04138402976425696
$ws.Range("Q5").Value = 30.93347286555778
$ws.Range("R5").Value = 278.40125579002
$ws.Range("S5").Value = 0.006066976641914681
$ws.Range("T5").Value = 0.006066976641914681
# Row 6
$ws.Range("I6").Value = 0.1466018818485066
$ws.Range("J6").Value = 0.1466018818485066
$ws.Range("O6").Value = 0.3297460182766552
$ws.Range("P6").Value = 0.3297460182766552
$ws.Range("S6").Value = 0.04834138681140972
$ws.Range("T6").Value = 0.04834138681140972
# Row 7
$ws.Range("I7").Value = 0.1466018818485066
$ws.Range("J7").Value = 0.1466018818485066
$ws.Range("O7").Value = 0.6288699519590879
$ws.Range("P7").Value = 0.6288699519590879
$ws.Range("S7").Value = 0.09219351839518225
$ws.Range("T7").Value = 0.09219351839518225
# Row 8
$ws.Range("I8").Value = 0.1672225590108659
$ws.Range("J8").Value = 0.1672225590108659
$ws.Range("M8").Value = 8.488196666666667
$ws.Range("N8").Value = 25.46459
$ws.Range("O8").Value = 0.04138402976425696
$ws.Range("P8").Value = 0.04138402976425696
$ws.Range("Q8").Value = 35.28450267109889
$ws.Range("R8").Value = 317.56052403989
$ws.Range("S8").Value = 0.006920343359360892
$ws.Range("T8").Value = 0.006920343359360892
# Row 9
$ws.Range("I9").Value = 0.1672225590108659
$ws.Range("J9").Value = 0.1672225590108659
$ws.Range("O9").Value = 0.3297460182766552
$ws.Range("P9").Value = 0.3297460182766552
$ws.Range("S9").Value = 0.05514097299986605
$ws.Range("T9").Value = 0.05514097299986605
# Row 10
$ws.Range("I10").Value = 0.1672225590108659
$ws.Range("J10").Value = 0.1672225590108659
$ws.Range("O10").Value = 0.6288699519590879
$ws.Range("P10").Value = 0.6288699519590879
$ws.Range("S10").Value = 0.105161242651639
$ws.Range("T10").Value = 0.105161242651639
